# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.683.95"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.939.48"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.80"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.21"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("E7").Value = "  -4.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -6.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.27"
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0842"
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.72"
$ws.Range("E13").Value = "  -5.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.417.98"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.40"
$ws.Range("E15").Value = "  -6.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.939.70"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.975"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.672.59"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.30"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -6.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.61"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.34"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.66"
$ws.Range("E25").Value = "  -5.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("E26").Value = "  -7.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.38"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.98"
$ws.Range("E32").Value = "  -6.25%  "
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.14"
$ws.Range("E34").Value = "  -7.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.72"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0423"
$ws.Range("E37").Value = "  -4.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.14"
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  -6.66%  "
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.61"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.59"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.100.22"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.22"
$ws.Range("E48").Value = "  -7.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.252.18"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.235"
$ws.Range("E50").Value = "  -5.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0320"
$ws.Range("E51").Value = "  -3.98%  "
